$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3748.5
$ws.Range("J32").Value = 2949.1667
$ws.Range("L32").Value = 2949.1667
$ws.Range("N32").Value = -3601.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3079.8333
$ws.Range("I111").Value = 3419.75
$ws.Range("J111").Value = 2400
$ws.Range("K111").Value = 10259.25
$ws.Range("L111").Value = 7200
$ws.Range("M111").Value = -7192.25
$ws.Range("N111").Value = -13334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 948.4194
$ws.Range("I125").Value = 676.8182
$ws.Range("J125").Value = 1612.3334
$ws.Range("K125").Value = 6091.3638
$ws.Range("L125").Value = 14511.0006
$ws.Range("M125").Value = -3631.3638
$ws.Range("N125").Value = -19431.0006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 26318760
$ws.Range("I137").Value = 1907.8462
$ws.Range("J137").Value = 40003524
$ws.Range("K137").Value = 5723.5386
$ws.Range("L137").Value = 120010572
$ws.Range("M137").Value = -3173.5386
$ws.Range("N137").Value = -120015672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15558.682
$ws.Range("I32").Value = 12962.863
$ws.Range("J32").Value = 72666.664
$ws.Range("K32").Value = 12962.863
$ws.Range("L32").Value = 72666.664
$ws.Range("M32").Value = -12675.863
$ws.Range("N32").Value = -73240.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9845.333000000001
$ws.Range("J37").Value = 11407.6
$ws.Range("L37").Value = 11407.6
$ws.Range("N37").Value = -11953.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14711600
$ws.Range("I74").Value = 20834612
$ws.Range("J74").Value = 16369
$ws.Range("K74").Value = 20834612
$ws.Range("L74").Value = 16369
$ws.Range("M74").Value = -20833738
$ws.Range("N74").Value = -18117

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14711600
$ws.Range("I77").Value = 20834612
$ws.Range("J77").Value = 16369
$ws.Range("K77").Value = 104173060
$ws.Range("L77").Value = 81845
$ws.Range("M77").Value = -104168692
$ws.Range("N77").Value = -90581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2458.3333
$ws.Range("I122").Value = 2600
$ws.Range("J122").Value = 2430
$ws.Range("K122").Value = 7800
$ws.Range("L122").Value = 7290
$ws.Range("M122").Value = -5350
$ws.Range("N122").Value = -12190

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2548.8572
$ws.Range("I31").Value = 1273.8182
$ws.Range("K31").Value = 1273.8182
$ws.Range("M31").Value = -978.8181999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2548.8572
$ws.Range("I34").Value = 1273.8182
$ws.Range("K34").Value = 1273.8182
$ws.Range("M34").Value = -1071.8182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10817.8
$ws.Range("J50").Value = 10817.8
$ws.Range("L50").Value = 10817.8
$ws.Range("N50").Value = -12067.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9668.799999999999
$ws.Range("J51").Value = 9886
$ws.Range("L51").Value = 9886
$ws.Range("N51").Value = -11358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 16233.333
$ws.Range("J59").Value = 16280
$ws.Range("L59").Value = 16280
$ws.Range("N59").Value = -18570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 10026.5
$ws.Range("J60").Value = 10026.5
$ws.Range("L60").Value = 10026.5
$ws.Range("N60").Value = -11048.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9668.799999999999
$ws.Range("J61").Value = 9886
$ws.Range("L61").Value = 9886
$ws.Range("N61").Value = -10582

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 14398
$ws.Range("J74").Value = 17469
$ws.Range("L74").Value = 17469
$ws.Range("N74").Value = -19217

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 14398
$ws.Range("J77").Value = 17469
$ws.Range("L77").Value = 52407
$ws.Range("N77").Value = -61143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 112536.22
$ws.Range("I99").Value = 1442.4
$ws.Range("J99").Value = 251403.5
$ws.Range("K99").Value = 1442.4
$ws.Range("L99").Value = 251403.5
$ws.Range("M99").Value = 55.59999999999991
$ws.Range("N99").Value = -254399.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 11601.2
$ws.Range("I122").Value = 15816
$ws.Range("J122").Value = 1766.6666
$ws.Range("K122").Value = 47448
$ws.Range("L122").Value = 5299.9998
$ws.Range("M122").Value = -44998
$ws.Range("N122").Value = -10199.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 112536.22
$ws.Range("I126").Value = 1442.4
$ws.Range("J126").Value = 251403.5
$ws.Range("K126").Value = 4327.200000000001
$ws.Range("L126").Value = 754210.5
$ws.Range("M126").Value = -1857.200000000001
$ws.Range("N126").Value = -759150.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15152128
$ws.Range("I5").Value = 620.8261
$ws.Range("J5").Value = 50000596
$ws.Range("K5").Value = 1862.4783
$ws.Range("L5").Value = 150001788
$ws.Range("M5").Value = -1750.4783
$ws.Range("N5").Value = -150002012

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 289.25
$ws.Range("I98").Value = 305
$ws.Range("J98").Value = 242
$ws.Range("K98").Value = 915
$ws.Range("L98").Value = 726
$ws.Range("M98").Value = 583
$ws.Range("N98").Value = -3722

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1221.3334
$ws.Range("J129").Value = 1602.2727
$ws.Range("L129").Value = 4806.8181
$ws.Range("N129").Value = -14806.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 15152128
$ws.Range("I135").Value = 620.8261
$ws.Range("J135").Value = 50000596
$ws.Range("K135").Value = 5587.4349
$ws.Range("L135").Value = 450005364
$ws.Range("M135").Value = -3052.4349
$ws.Range("N135").Value = -450010434

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 368.3
$ws.Range("I107").Value = 298.1111
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 298.1111
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1621.8889
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3837.5652
$ws.Range("I122").Value = 4077.1052
$ws.Range("J122").Value = 2699.75
$ws.Range("K122").Value = 12231.3156
$ws.Range("L122").Value = 8099.25
$ws.Range("M122").Value = -9781.3156
$ws.Range("N122").Value = -12999.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1526.8572
$ws.Range("I61").Value = 1314.6666
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 1314.6666
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -1112.6666
$ws.Range("N61").Value = -3204

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1526.8572
$ws.Range("I113").Value = 1314.6666
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 1314.6666
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = 855.3334
$ws.Range("N113").Value = -7140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2574.0637
$ws.Range("I132").Value = 2923.3547
$ws.Range("K132").Value = 8770.0641
$ws.Range("M132").Value = -6240.0641
